$wb = $excel.ActiveWorkbook

# Trade #32 closed/new row added to both the "All Trades" log and the
# per-strategy "base_strategy" sheet (the workbook keeps a duplicate of
# the trade log on the strategy-specific tab).
$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 33

    $ws.Cells.Item($row, 1).Value = 32                     # A - Trade #
    $ws.Cells.Item($row, 2).Value = "'2026-02-16"          # B - Date (text, not a date serial)
    $ws.Cells.Item($row, 3).Value = "'22:55:17"            # C - Time (text, not a time serial)
    $ws.Cells.Item($row, 4).Value = "base_strategy"        # D - Strategy
    $ws.Cells.Item($row, 5).Value = "DOWN"                 # E - Side
    $ws.Cells.Item($row, 6).Value = 49.999998              # F - Entry Price
    # G - Exit Price stays blank (trade still OPEN)
    $ws.Cells.Item($row, 8).Value = "OPEN"                 # H - Status
    $ws.Cells.Item($row, 9).Value = 0                      # I - P&L %
    $ws.Cells.Item($row, 10).Value = 0                     # J - P&L $
    $ws.Cells.Item($row, 11).Value = 100                   # K - Capital After
    $ws.Cells.Item($row, 12).Value = 0                     # L - Entry Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0                     # M - Exit Slippage (bps)
    $ws.Cells.Item($row, 14).Value = 0.6                   # N - Confidence
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"  # O - Entry Reason
    # P - Exit Reason stays blank (trade still OPEN)
    $ws.Cells.Item($row, 17).Value = 0                     # Q - Duration (min)
}
